$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the coordinate values in Q2 and R2
$ws.Range("Q2").Value = 692392
$ws.Range("R2").Value = 6610760

# Clear the "Starttid" (Z2) and "Sluttid" (AB2) time values entirely
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
